# 8.10.2.2 "Insurance companies' financial indicators" sheet update:
#  - retitle the sheet heading (Kyrgyz text, same shared-string slot)
#  - append a new "2023" data column (Q) mirroring the existing "2022"
#    column (P): header year, reporting-company count, premium amount
#  - tidy up the stale selection left over from editing (was E9, an
#    empty cell below the table; move it back to A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Retitle the report header (row 1, column A) ---------------------
$ws.Range("A1").Value = "8.10.2.2 Камсыздандыруу компаниялардын финансылык көрсөткүчтөрү"

# --- New 2023 column: copy formatting from the 2022 column (P) then
#     fill in the 2023 figures ------------------------------------------------
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2023

$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 16

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 3031.4

# --- Reset selection to A1 (was pointing at a stale cell) -------------
$ws.Range("A1").Select()
